$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-24 Monday" "2025-02-25 Tuesday"

Replace-Text "605×7=" "420×5="
Replace-Text "358×8=" "126×8="
Replace-Text "408×3=" "741×7="
Replace-Text "150×5=" "788×8="
Replace-Text "232×8=" "722×3="
Replace-Text "625×8=" "675×3="
Replace-Text "539×5=" "537×5="
Replace-Text "676×3=" "572×9="
Replace-Text "811×8=" "330×6="
Replace-Text "715×8=" "207×5="
Replace-Text "340×6=" "141×7="
Replace-Text "902×9=" "237×9="
Replace-Text "722×2=" "832×5="
Replace-Text "521×4=" "326×9="
Replace-Text "903×5=" "779×7="
Replace-Text "137×2=" "223×3="
Replace-Text "148×2=" "511×6="
Replace-Text "729×6=" "691×3="
Replace-Text "748×9=" "545×9="
Replace-Text "232×3=" "394×8="
Replace-Text "449×7=" "181×3="
Replace-Text "162×4=" "585×5="
Replace-Text "722×6=" "318×8="
Replace-Text "803×9=" "854×3="
Replace-Text "238×6=" "482×9="
